# Adds a new "item_num" column (col O) ahead of the existing "comment"
# column (which shifts from O to P) on every sheet that carries the
# common Measures-style header (Measures, ID, Dems, Dates, NewVars).
# Also restores the data entered for item_num on the Measures sheet,
# and leaves NewVars as the active sheet/tab, matching the saved view
# state from the author's last edit.

$wb = $excel.ActiveWorkbook

# --- Measures --------------------------------------------------------
$ws = $wb.Worksheets.Item("Measures")
$ws.Columns.Item(15).Insert()
$ws.Range("O1").Value = "item_num"
$ws.Range("O2").Value = 1
$ws.Range("O3").Value = 1
$ws.Range("O4").Value = 1
$ws.Range("O5").Value = 4
$ws.Range("O6").Value = 4
$ws.Activate()
$ws.Range("O7").Select()

# --- ID ----------------------------------------------------------------
$ws = $wb.Worksheets.Item("ID")
$ws.Columns.Item(15).Insert()
$ws.Range("O1").Value = "item_num"
$ws.Activate()
$ws.Range("O2").Select()

# --- Dems ----------------------------------------------------------------
$ws = $wb.Worksheets.Item("Dems")
$ws.Columns.Item(15).Insert()
$ws.Range("O1").Value = "item_num"
$ws.Activate()
$ws.Range("O2").Select()

# --- Dates ----------------------------------------------------------------
$ws = $wb.Worksheets.Item("Dates")
$ws.Columns.Item(15).Insert()
$ws.Range("O1").Value = "item_num"
$ws.Activate()
$ws.Range("O2").Select()

# --- NewVars (left as the active sheet/tab, matching saved state) -------
$ws = $wb.Worksheets.Item("NewVars")
$ws.Columns.Item(15).Insert()
$ws.Range("O1").Value = "item_num"
$ws.Activate()
$ws.Range("G10").Select()
